$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in unique "this is test8" .. "this is test18" values for rows 10..20 (column A)
for ($i = 0; $i -le 10; $i++) {
    $row = 10 + $i
    $n = 8 + $i
    $ws.Cells.Item($row, 1).Value = "this is test$n"
}

# Move the active selection from A9 to A2
$ws.Range("A2").Select()
